# Scheduled-runner style update of market price / profit figures across
# the per-crafting-class Leve tables (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Each block below refreshes currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ)
# and the derived LeveProfit(NQ/HQ) columns (H..N) for the affected rows.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 103.083336
$ws.Range("I9").Value = 103.083336
$ws.Range("K9").Value = 103.083336
$ws.Range("M9").Value = 65.916664

$ws.Range("H33").Value = 222.06667
$ws.Range("I33").Value = 202.21428
$ws.Range("K33").Value = 202.21428
$ws.Range("M33").Value = 26.78572

$ws.Range("H58").Value = 2028.5555
$ws.Range("J58").Value = 2531.1667
$ws.Range("L58").Value = 7593.500100000001
$ws.Range("N58").Value = -7893.500100000001

$ws.Range("H82").Value = 1209.4
$ws.Range("I82").Value = 1209.4
$ws.Range("K82").Value = 3628.2
$ws.Range("M82").Value = -3222.2

$ws.Range("H85").Value = 1209.4
$ws.Range("I85").Value = 1209.4
$ws.Range("K85").Value = 3628.2
$ws.Range("M85").Value = -2224.2

$ws.Range("H112").Value = 3036.8333
$ws.Range("J112").Value = 2716
$ws.Range("L112").Value = 8148
$ws.Range("N112").Value = -10364

$ws.Range("H138").Value = 5486.7
$ws.Range("I138").Value = 5692.4287
$ws.Range("J138").Value = 5375.923
$ws.Range("K138").Value = 17077.2861
$ws.Range("L138").Value = 16127.769
$ws.Range("M138").Value = -11937.2861
$ws.Range("N138").Value = -26407.769

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H60").Value = 5000
$ws.Range("I60").Value = 5000
$ws.Range("K60").Value = 5000
$ws.Range("M60").Value = -4267

$ws.Range("H61").Value = 2401.0908
$ws.Range("I61").Value = 2401.0908
$ws.Range("K61").Value = 2401.0908
$ws.Range("M61").Value = -2189.0908

$ws.Range("H74").Value = 4657.132
$ws.Range("I74").Value = 2453.0908
$ws.Range("K74").Value = 2453.0908
$ws.Range("M74").Value = -1579.0908

$ws.Range("H77").Value = 4657.132
$ws.Range("I77").Value = 2453.0908
$ws.Range("K77").Value = 12265.454
$ws.Range("M77").Value = -7897.454

$ws.Range("H110").Value = 1225
$ws.Range("I110").Value = 1225
$ws.Range("K110").Value = 1225
$ws.Range("M110").Value = 820

$ws.Range("H122").Value = 2274.1292
$ws.Range("I122").Value = 1868.36
$ws.Range("K122").Value = 5605.08
$ws.Range("M122").Value = -3155.08

$ws.Range("H136").Value = 2401.0908
$ws.Range("I136").Value = 2401.0908
$ws.Range("K136").Value = 7203.2724
$ws.Range("M136").Value = -4653.2724

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H56").Value = 20000
$ws.Range("J56").Value = 20000
$ws.Range("L56").Value = 20000
$ws.Range("N56").Value = -21478

$ws.Range("H94").Value = 5178.8
$ws.Range("I94").Value = 5059.4165
$ws.Range("K94").Value = 5059.4165
$ws.Range("M94").Value = -4608.4165

$ws.Range("H107").Value = 1224.5714
$ws.Range("I107").Value = 1224.5714
$ws.Range("K107").Value = 1224.5714
$ws.Range("M107").Value = 695.4286

$ws.Range("H134").Value = 896.8276
$ws.Range("I134").Value = 819.96
$ws.Range("J134").Value = 1377.25
$ws.Range("K134").Value = 2459.88
$ws.Range("L134").Value = 4131.75
$ws.Range("M134").Value = 75.11999999999989
$ws.Range("N134").Value = -9201.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H29").Value = 5250
$ws.Range("I29").Value = 5500
$ws.Range("K29").Value = 5500
$ws.Range("M29").Value = -5207

$ws.Range("H31").Value = 3759.9
$ws.Range("I31").Value = 3759.9
$ws.Range("K31").Value = 3759.9
$ws.Range("M31").Value = -3464.9

$ws.Range("H34").Value = 3759.9
$ws.Range("I34").Value = 3759.9
$ws.Range("K34").Value = 3759.9
$ws.Range("M34").Value = -3557.9

$ws.Range("H43").Value = 17600
$ws.Range("J43").Value = 17600
$ws.Range("L43").Value = 17600
$ws.Range("N43").Value = -17968

$ws.Range("H101").Value = 17600
$ws.Range("J101").Value = 17600
$ws.Range("L101").Value = 17600
$ws.Range("N101").Value = -24090

$ws.Range("H110").Value = 88000
$ws.Range("J110").Value = 88000
$ws.Range("L110").Value = 88000
$ws.Range("N110").Value = -96180

$ws.Range("H134").Value = 2740.3
$ws.Range("I134").Value = 1901.5
$ws.Range("K134").Value = 5704.5
$ws.Range("M134").Value = -3169.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 399.42105
$ws.Range("J12").Value = 102.84615
$ws.Range("L12").Value = 308.53845
$ws.Range("N12").Value = -654.53845

$ws.Range("H34").Value = 1807.125
$ws.Range("I34").Value = 1332.5
$ws.Range("J34").Value = 1965.3334
$ws.Range("K34").Value = 3997.5
$ws.Range("L34").Value = 5896.0002
$ws.Range("M34").Value = -3913.5
$ws.Range("N34").Value = -6064.0002

$ws.Range("H37").Value = 41362.727
$ws.Range("J37").Value = 41362.727
$ws.Range("L37").Value = 124088.181
$ws.Range("N37").Value = -124312.181

$ws.Range("H39").Value = 116154.445
$ws.Range("I39").Value = 151098
$ws.Range("J39").Value = 72475
$ws.Range("K39").Value = 453294
$ws.Range("L39").Value = 217425
$ws.Range("M39").Value = -453000
$ws.Range("N39").Value = -218013

$ws.Range("H69").Value = 2433.125
$ws.Range("I69").Value = 1878
$ws.Range("K69").Value = 5634
$ws.Range("M69").Value = -4823

$ws.Range("H72").Value = 2433.125
$ws.Range("I72").Value = 1878
$ws.Range("K72").Value = 16902
$ws.Range("M72").Value = -12846

$ws.Range("H97").Value = 568.25
$ws.Range("I97").Value = 424.66666
$ws.Range("K97").Value = 1273.99998
$ws.Range("M97").Value = -777.9999800000001

$ws.Range("H134").Value = 7229.1113
$ws.Range("I134").Value = 3960.75
$ws.Range("J134").Value = 13765.833
$ws.Range("K134").Value = 11882.25
$ws.Range("L134").Value = 41297.499
$ws.Range("M134").Value = -6812.25
$ws.Range("N134").Value = -51437.499

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 1423.8
$ws.Range("I31").Value = 1423.8
$ws.Range("K31").Value = 1423.8
$ws.Range("M31").Value = -1131.8

$ws.Range("H37").Value = 1423.8
$ws.Range("I37").Value = 1423.8
$ws.Range("K37").Value = 1423.8
$ws.Range("M37").Value = -1146.8

$ws.Range("H70").Value = 25402.834
$ws.Range("I70").Value = 28996.223
$ws.Range("J70").Value = 14622.667
$ws.Range("K70").Value = 28996.223
$ws.Range("L70").Value = 14622.667
$ws.Range("M70").Value = -28726.223
$ws.Range("N70").Value = -15162.667

$ws.Range("H73").Value = 25402.834
$ws.Range("I73").Value = 28996.223
$ws.Range("J73").Value = 14622.667
$ws.Range("K73").Value = 28996.223
$ws.Range("L73").Value = 14622.667
$ws.Range("M73").Value = -28060.223
$ws.Range("N73").Value = -16494.667

$ws.Range("H80").Value = 4569.3335
$ws.Range("I80").Value = 4602.5
$ws.Range("J80").Value = 4503
$ws.Range("K80").Value = 4602.5
$ws.Range("L80").Value = 4503
$ws.Range("M80").Value = -3604.5
$ws.Range("N80").Value = -6499

$ws.Range("H83").Value = 4569.3335
$ws.Range("I83").Value = 4602.5
$ws.Range("J83").Value = 4503
$ws.Range("K83").Value = 23012.5
$ws.Range("L83").Value = 22515
$ws.Range("M83").Value = -18020.5
$ws.Range("N83").Value = -32499

$ws.Range("H122").Value = 2367.2068
$ws.Range("I122").Value = 1679.5454
$ws.Range("J122").Value = 2787.4443
$ws.Range("K122").Value = 5038.6362
$ws.Range("L122").Value = 8362.332900000001
$ws.Range("M122").Value = -2588.6362
$ws.Range("N122").Value = -13262.3329

$ws.Range("H132").Value = 11797.857
$ws.Range("I132").Value = 11906.061
$ws.Range("J132").Value = 10012.5
$ws.Range("K132").Value = 35718.183
$ws.Range("L132").Value = 30037.5
$ws.Range("M132").Value = -33188.183
$ws.Range("N132").Value = -35097.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 34700000
$ws.Range("I2").Value = 2050002
$ws.Range("K2").Value = 2050002
$ws.Range("M2").Value = -2049890

$ws.Range("H132").Value = 4677.7144
$ws.Range("I132").Value = 3753.3333
$ws.Range("J132").Value = 5371
$ws.Range("K132").Value = 11259.9999
$ws.Range("L132").Value = 16113
$ws.Range("M132").Value = -8729.999899999999
$ws.Range("N132").Value = -21173

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 100000
$ws.Range("J3").Value = 100000
$ws.Range("L3").Value = 100000
$ws.Range("N3").Value = -100228

$ws.Range("H22").Value = 999
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

$ws.Range("H132").Value = 4265.8887
$ws.Range("I132").Value = 4032.1667
$ws.Range("K132").Value = 12096.5001
$ws.Range("M132").Value = -9566.500100000001
